$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title / info block (rows 1-8, merged A:H) ----
$ws.Range("A1").Value = "Clase de formacion"
$ws.Range("A6").Value = "Instructor: Jhon Becerra"
$ws.Range("A7").Value = "Clase Formacion: Desarrollo Web"
$ws.Range("A8").Value = "Fecha: 2024-10-15 12:09:17"

# ---- Row 11 (attendee 1) ----
$ws.Range("B11").Value = "'0987654321"
$ws.Range("C11").Value = "Enernesto perez"
$ws.Range("D11").Value = "'111111111"
$ws.Range("E11").Value = "ernesto@gmail.com"
$ws.Range("G11").Value = "Santander - Bucaramanga - El Llano"

# ---- Row 12 (attendee 2) ----
$ws.Range("B12").Value = "'1234567890"
$ws.Range("C12").Value = "Carlos ALberto Torrez"
$ws.Range("D12").Value = "'8912381297"
$ws.Range("E12").Value = "calberto@gmail.com"
$ws.Range("F12").Value = "Masculino"
$ws.Range("G12").Value = "Santander - Bucaramanga - La Esperanza"
$ws.Range("H12").Value = 0
$ws.Range("H12").Interior.ColorIndex = 35

# ---- Row 13 (attendee 3) ----
$ws.Range("B13").Value = "'99999999"
$ws.Range("C13").Value = "Pedro albaro Quinteroo"
$ws.Range("D13").Value = "'112222221"
$ws.Range("E13").Value = "pedro@gmail.com"
$ws.Range("G13").Value = "Santander - Bucaramanga - El Llano"
$ws.Range("H13").Value = 5
$ws.Range("H13").Interior.ColorIndex = 3

# ---- Row 14 (new attendee 4) ----
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "'00000000"
$ws.Range("C14").Value = "ejero alfredo torrez alcanso"
$ws.Range("D14").Value = "'098098098"
$ws.Range("E14").Value = "algo@gmail.com"
$ws.Range("F14").Value = "Masculino"
$ws.Range("G14").Value = "Santander - Bucaramanga - El Llano"
$ws.Range("H14").Value = 5
$ws.Range("H14").Interior.ColorIndex = 3

# ---- Column width tweaks (best-fit widths recomputed for the new data) ----
$ws.Columns.Item(3).ColumnWidth = 24.75
$ws.Columns.Item(5).ColumnWidth = 18.75
$ws.Columns.Item(6).ColumnWidth = 9.25
